# record model add apply member field
#
# Both sheets gain a new "申請人數" (number of applicants) column, inserted
# just before the existing "核准人數" column. That shifts every column from
# the (old) F position onward one slot to the right, growing each sheet's
# used range by one column (K -> L).

$wb = $excel.ActiveWorkbook

# ---- Sheet "總數" (template header sheet) ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Columns("F:F").Insert()
$ws1.Range("F4").Value = "申請人數"

# ---- Sheet "區域月份統計" (monthly/region stats sheet) ----
$ws2 = $wb.Worksheets.Item(2)
$ws2.Columns("F:F").Insert()
$ws2.Range("F2").Value = "申請人數"

# Restore a sane selection/cursor position on each sheet.
$ws2.Range("H9").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("A1").Select() | Out-Null
